$d = $word.ActiveDocument

# ----------------------------------------------------------------------
# Step 1: Merge the split "arrow" runs ("->" + " Word") into single runs.
# Word's COM Find/Replace (Replace = wdReplaceAll) naturally merges a run
# that had its text replaced with adjacent runs sharing identical
# formatting, which is exactly the cleanup performed in the diff.
# Running a (no-op content-wise) replace-all on the arrow glyph triggers
# that merge everywhere the adjacent run has matching rPr, and correctly
# leaves the "getName ... -> Name" pair untouched (their rPr differs
# because of the surrounding proofErr/gramEnd markers breaking the run
# boundary there), matching the diff exactly.
# ----------------------------------------------------------------------
$arrow = [char]0x2192
$null = $d.Content.Find.Execute($arrow, $true, $false, $false, $false, $false, `
                                 $true, 1, $false, $arrow, 2)

# ----------------------------------------------------------------------
# Step 2: Merge the 3-way split "{inv: supplies>=0}" runs in the TAD
# Weapon class invariant row the same way.
# ----------------------------------------------------------------------
$null = $d.Content.Find.Execute("supplies>=0", $true, $false, $false, $false, $false, `
                                 $true, 1, $false, "supplies>=0", 2)

# ----------------------------------------------------------------------
# Step 3: Insert the new "isEmpty ... Supplies -> booleano" bullet point
# right after the "useWeapon" bullet in the TAD Weapon operations cell,
# carrying the _GoBack bookmark that used to sit at the end of the
# document.
# ----------------------------------------------------------------------
$paras = $d.Paragraphs
$targetIdx = -1
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -like "useWeapon*Weapon x Supplies*") {
        $targetIdx = $i
    }
}
$useWeaponPara = $paras.Item($targetIdx)
$useWeaponPara.Range.InsertParagraphAfter()

$paras = $d.Paragraphs
$newPara = $paras.Item($targetIdx + 1)

$newParaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>isEmpty</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">                           Supplies </w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:lang w:val="en-US"/></w:rPr><w:t>→</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> booleano</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newPara.Range.InsertXML($newParaXml)

# ----------------------------------------------------------------------
# Step 4: Remove the _GoBack bookmark from its old home (the empty
# paragraph right before the final sectPr), while leaving that paragraph
# itself (and its formatting) intact. We do this by overwriting its
# content with throw-away marker text via InsertXML (which - because the
# inserted content is itself a full paragraph - pushes the original
# paragraph mark/properties, now stripped of the bookmark it held, back
# to the end) and then deleting the leftover marker paragraph.
# ----------------------------------------------------------------------
$paras = $d.Paragraphs
$lastPara = $paras.Item($paras.Count)

$markerXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>ZZMARKERZZ</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$lastPara.Range.InsertXML($markerXml)

$findRange = $d.Content
$found = $findRange.Find.Execute("ZZMARKERZZ", $true, $false, $false, $false, $false, `
                                  $true, 1, $false, "", 0)
if ($found) {
    $delRange = $d.Range($findRange.Start, $findRange.End + 1)
    $delRange.Delete()
}
